$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 200.4
$ws.Range("I9").Value = 166.66667
$ws.Range("J9").Value = 251
$ws.Range("K9").Value = 166.66667
$ws.Range("L9").Value = 251
$ws.Range("M9").Value = 2.333329999999989
$ws.Range("N9").Value = -589
$ws.Range("H17").Value = 1800
$ws.Range("J17").Value = 1800
$ws.Range("L17").Value = 5400
$ws.Range("N17").Value = -5736
$ws.Range("H76").Value = 6783.25
$ws.Range("I76").Value = 4633
$ws.Range("J76").Value = 7500
$ws.Range("K76").Value = 4633
$ws.Range("L76").Value = 7500
$ws.Range("M76").Value = -4318
$ws.Range("N76").Value = -8130
$ws.Range("H79").Value = 6783.25
$ws.Range("I79").Value = 4633
$ws.Range("J79").Value = 7500
$ws.Range("K79").Value = 4633
$ws.Range("L79").Value = 7500
$ws.Range("M79").Value = -3541
$ws.Range("N79").Value = -9684
$ws.Range("H138").Value = 2247.785
$ws.Range("I138").Value = 822.5833
$ws.Range("J138").Value = 2869.691
$ws.Range("K138").Value = 2467.7499
$ws.Range("L138").Value = 8609.073
$ws.Range("M138").Value = 2672.2501
$ws.Range("N138").Value = -18889.073
$ws.Range("H141").Value = 9390.125
$ws.Range("I141").Value = 10161.571
$ws.Range("K141").Value = 30484.713
$ws.Range("M141").Value = -25304.713
# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 634.16
$ws.Range("I2").Value = 574.6667
$ws.Range("K2").Value = 574.6667
$ws.Range("M2").Value = -461.6667
$ws.Range("H45").Value = 21741574
$ws.Range("I45").Value = 29413590
$ws.Range("J45").Value = 4196.6665
$ws.Range("K45").Value = 29413590
$ws.Range("L45").Value = 4196.6665
$ws.Range("M45").Value = -29413213
$ws.Range("N45").Value = -4950.6665
$ws.Range("H92").Value = 70499.60000000001
$ws.Range("J92").Value = 70499.60000000001
$ws.Range("L92").Value = 70499.60000000001
$ws.Range("N92").Value = -75491.60000000001
$ws.Range("H94").Value = 47990
$ws.Range("J94").Value = 47990
$ws.Range("L94").Value = 47990
$ws.Range("N94").Value = -49792
$ws.Range("H116").Value = 634.16
$ws.Range("I116").Value = 574.6667
$ws.Range("K116").Value = 574.6667
$ws.Range("M116").Value = 1719.3333
# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 634.16
$ws.Range("I3").Value = 574.6667
$ws.Range("K3").Value = 574.6667
$ws.Range("M3").Value = -460.6667
$ws.Range("H54").Value = 24999.5
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H76").Value = 85000.5
$ws.Range("J76").Value = 85000.5
$ws.Range("L76").Value = 85000.5
$ws.Range("N76").Value = -85630.5
$ws.Range("H79").Value = 85000.5
$ws.Range("J79").Value = 85000.5
$ws.Range("L79").Value = 85000.5
$ws.Range("N79").Value = -87184.5
# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 124.5
$ws.Range("I26").Value = 124.5
$ws.Range("K26").Value = 373.5
$ws.Range("M26").Value = -85.5
$ws.Range("H46").Value = 2167.5715
$ws.Range("J46").Value = 2361.6667
$ws.Range("L46").Value = 7085.000100000001
$ws.Range("N46").Value = -7267.000100000001
$ws.Range("H131").Value = 6294.7
$ws.Range("J131").Value = 6294.7
$ws.Range("L131").Value = 18884.1
$ws.Range("N131").Value = -28964.1
# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 7650
$ws.Range("H57").Value = 16809.6
$ws.Range("I57").Value = 4055
$ws.Range("J57").Value = 19998.25
$ws.Range("K57").Value = 4055
$ws.Range("L57").Value = 19998.25
$ws.Range("M57").Value = -3235
$ws.Range("N57").Value = -21638.25
$ws.Range("H62").Value = 115000
$ws.Range("J62").Value = 115000
$ws.Range("L62").Value = 115000
$ws.Range("N62").Value = -116372
$ws.Range("H65").Value = 115000
$ws.Range("J65").Value = 115000
$ws.Range("L65").Value = 345000
$ws.Range("N65").Value = -351864
$ws.Range("H70").Value = 9353
$ws.Range("I70").Value = 10983.571
$ws.Range("K70").Value = 10983.571
$ws.Range("M70").Value = -10713.571
$ws.Range("H73").Value = 9353
$ws.Range("I73").Value = 10983.571
$ws.Range("K73").Value = 10983.571
$ws.Range("M73").Value = -10047.571
$ws.Range("H80").Value = 4247.909
$ws.Range("I80").Value = 3342.25
$ws.Range("K80").Value = 3342.25
$ws.Range("M80").Value = -2344.25
$ws.Range("H83").Value = 4247.909
$ws.Range("I83").Value = 3342.25
$ws.Range("K83").Value = 16711.25
$ws.Range("M83").Value = -11719.25
$ws.Range("H93").Value = 64999
$ws.Range("J93").Value = 64999
$ws.Range("L93").Value = 64999
$ws.Range("N93").Value = -68743
$ws.Range("H123").Value = 52000
$ws.Range("J123").Value = 52000
$ws.Range("L123").Value = 52000
$ws.Range("N123").Value = -56900
$ws.Range("H136").Value = 10681.5
$ws.Range("J136").Value = 10681.5
$ws.Range("L136").Value = 32044.5
$ws.Range("N136").Value = -37144.5
# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1424.9565
$ws.Range("J22").Value = 1638.625
$ws.Range("L22").Value = 1638.625
$ws.Range("N22").Value = -2228.625
$ws.Range("H27").Value = 1424.9565
$ws.Range("J27").Value = 1638.625
$ws.Range("L27").Value = 1638.625
$ws.Range("N27").Value = -1852.625
$ws.Range("H46").Value = 3081.875
$ws.Range("I46").Value = 2203.75
$ws.Range("J46").Value = 4838.125
$ws.Range("K46").Value = 2203.75
$ws.Range("L46").Value = 4838.125
$ws.Range("M46").Value = -2015.75
$ws.Range("N46").Value = -5214.125
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H76").Value = 30000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 30000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 30000
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -30676
$ws.Range("H79").Value = 30000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 30000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 30000
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -32340
$ws.Range("H115").Value = 65000
$ws.Range("J115").Value = 65000
$ws.Range("L115").Value = 65000
$ws.Range("N115").Value = -67350
# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2256.4
$ws.Range("I122").Value = 2383.25
$ws.Range("J122").Value = 1749
$ws.Range("K122").Value = 7149.75
$ws.Range("L122").Value = 5247
$ws.Range("M122").Value = -4699.75
$ws.Range("N122").Value = -10147
